# Emmersive localization workbook update
# - Updates the "Provider Nickname" string (color=red -> i tag) and shrinks its row height
# - Adds four new localization rows (em_ui_config_reset, em_ui_config_open,
#   em_ui_config_changed, em_character_data_statuses)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General")

# --- Row 17: Provider Nickname string changed from <color=red> to <i> tag ---
$ws.Range("C17").Value = "プロバイダーニックネーム (<i>変更後にパラメータを再編集する必要がある</i>)"
$ws.Range("D17").Value = "Provider Nickname (<i>Changes will reset params</i>) "
$ws.Rows.Item(17).RowHeight = 23.25

# --- Row 82: new "Elin with AI configuration reset" strings ---
$ws.Range("A82").Value = "em_ui_config_reset"
$ws.Range("C82").Value = "「Elin with AI」の設定が{0}にリセットされました"
$ws.Range("D82").Value = "Elin with AI configuration has been resetted to {0}"

# --- Row 83: new "open mod configuration" strings ---
$ws.Range("A83").Value = "em_ui_config_open"
$ws.Range("C83").Value = "Mod設定を開く"
$ws.Range("D83").Value = "Edit Configuration"

# --- Row 84: new "configuration changed" strings ---
$ws.Range("A84").Value = "em_ui_config_changed"
$ws.Range("C84").Value = "Emmersive configuration changed"
$ws.Range("D84").Value = "Emmersive configuration changed"

# --- Row 85: new "character statuses" strings ---
$ws.Range("A85").Value = "em_character_data_statuses"
$ws.Range("C85").Value = "ステータス"
$ws.Range("D85").Value = "Statuses"

# --- Update view state to match where the author left the selection ---
$ws.Range("D85").Select()
$excel.ActiveWindow.ScrollRow = 71

Write-Output "Applied Emmersive localization updates"
